# Update "想去人数" (F column) counts that changed between crawler runs.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 131
$wsExhibit.Range("F3").Value = 2152
$wsExhibit.Range("F5").Value = 11334
$wsExhibit.Range("F6").Value = 202
$wsExhibit.Range("F9").Value = 11273
$wsExhibit.Range("F16").Value = 3469
$wsExhibit.Range("F17").Value = 174

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 131
$wsAll.Range("F3").Value = 2152
$wsAll.Range("F7").Value = 11334
$wsAll.Range("F8").Value = 202
$wsAll.Range("F11").Value = 11273
$wsAll.Range("F18").Value = 3469
$wsAll.Range("F19").Value = 174
